# Commit: "added author, career, country, quotes and subject to database"
#
# The "Quote.csv" staging sheet (a raw CSV-import mirror of the Quotes
# sheet) was missing trailing Subject2_ID/Subject3_ID numbers (columns F/G)
# for a batch of rows -- those gaps get back-filled with 0 here. Also
# restores the two UI selection rectangles that moved as part of the same
# editing session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Quote.csv: fill the previously-blank F (Subject2_ID) / G (Subject3_ID)
#    cells for rows 9-101 with 0.
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("Quote.csv")

$rowsF = @(15,16,23,26,30,31,34,38,41,46,47,51,55,60,62,64,65,70,71,73,77,80,82,84,88,90,92,95,96,99,100)
$rowsG = @(9,10,14,15,16,17,20,22,23,25,26,28,30,31,33,34,35,36,38,40,41,45,46,47,48,49,50,51,52,53,55,56,57,59,60,62,63,64,65,68,69,70,71,72,73,75,76,77,79,80,82,84,86,88,89,90,92,93,94,95,96,97,99,100,101)

foreach ($r in $rowsF) {
    $ws7.Cells.Item($r, 6).Value = 0
}
foreach ($r in $rowsG) {
    $ws7.Cells.Item($r, 7).Value = 0
}

# ---------------------------------------------------------------------
# 2) Restore the UI selection rectangles.
#    Order matters: touch "Quotes" first, then finish on "Quote.csv" so
#    that sheet ends up the active tab again (as in the source file).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Quotes")
[void]$ws2.Range("F1:J1048576").Select()

[void]$ws7.Range("M6").Select()
